# Total of 6 TCs (LoginTest.java):
# LoginWithAnalyzer, LoginWithIndividual, LoginWithUser,
# LoginWithSchemaManager, LoginWithUserManager, LoginWithSuperUser
#
# The "TestData" sheet already had columns for Analyzer / Individual /
# SchemaManager (B/C/D). This adds three more test-data columns
# (E/F/G) for the new User, UserManager and SuperUser test cases:
#   row 8  (Tenant)   -> demo / demo / demo
#   row 9  (Username) -> AbdelsalamUser / AbdelsalamUserManager / AbdelsalamSuper
#   row 10 (Password) -> AbdelsalamUser1 / AbdelsalamUserManager1 / AbdelsalamSuper1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 8 - Tenant
$ws.Range("E8").Value = "demo"
$ws.Range("F8").Value = "demo"
$ws.Range("G8").Value = "demo"

# Row 9 - Username
$ws.Range("E9").Value = "AbdelsalamUser"
$ws.Range("F9").Value = "AbdelsalamUserManager"
$ws.Range("G9").Value = "AbdelsalamSuper"

# Row 10 - Password
$ws.Range("E10").Value = "AbdelsalamUser1"
$ws.Range("F10").Value = "AbdelsalamUserManager1"
$ws.Range("G10").Value = "AbdelsalamSuper1"
